$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scheduled-task re-run rewrote the timestamp in row 16 with extra
# floating point precision (same instant, just re-serialized by Excel).
$ws.Cells.Item(16, 1).Value2 = 45865.70864775463

# Append the new row 17 captured by this run of the scheduled task.
$ws.Cells.Item(17, 1).Value2 = 45865.75030394648
$ws.Cells.Item(17, 1).NumberFormat = $ws.Cells.Item(16, 1).NumberFormat

$ws.Cells.Item(17, 2).Value = 2025
$ws.Cells.Item(17, 3).Value = 30
$ws.Cells.Item(17, 4).Value = 17.61
$ws.Cells.Item(17, 5).Value = 79.3
$ws.Cells.Item(17, 6).Value = 17.58
$ws.Cells.Item(17, 7).Value = 7.77
$ws.Cells.Item(17, 8).Value = "ESE"
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = "18:00:26"
